# Update experiment timestamps on the "Experiments" sheet (E2:E9)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Experiments")

$ws.Range("E2").Value = "07/06/2022 21:13:44"
$ws.Range("E3").Value = "07/06/2022 21:13:45"
$ws.Range("E4").Value = "07/06/2022 21:13:46"
$ws.Range("E5").Value = "07/06/2022 21:13:48"
$ws.Range("E6").Value = "07/06/2022 21:13:59"
$ws.Range("E7").Value = "07/06/2022 21:14:05"
$ws.Range("E8").Value = "07/06/2022 21:15:31"
$ws.Range("E9").Value = "07/06/2022 21:16:25"
